$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").NumberFormat = "@"
$ws.Range("B2").Value = "253"
$ws.Range("B2").Style = "Normal"
$ws.Range("C2").Value = "Yes, community sites are community-based physical locations where services and resources are provided to fulfill the needs of the society members."
$ws.Range("D2").Value = "No, Soleo Health is a healthcare provider company focused on providing specialty infusion services, and it does not engage in lobbying or advocacy efforts to influence state or local policies."
$ws.Range("E2").Value = "No, ,soleo Health does not provide engagement opportunity with leadership."
$ws.Range("F2").Value = "No, Soleo Health does not provide support for clinical trial recruitment. Soleo Health is a specialty pharmacy focusing on complex disease management and medication therapies."
$ws.Range("G2").Value = "No, Soleo Health does not provide engagement opportunity with payors. Soleo Health primarily focuses on providing specialty pharmacy and infusion services to patients, rather than direct engagement with payors."
$ws.Range("H2").Value = "No, Soleo Health does not include area experts on its board., The composition of the board does not indicate the presence of area experts specializing in specific fields."
$ws.Range("I2").Value = "Yes, Soleo Health is involved in therapeutic research collaborations. Soleo Health partners with various stakeholders in the healthcare industry to advance research and development in therapeutic treatments."
$ws.Range("J2").Value = "No, Soleo Health does not include top therapeutic area experts on its board. There is no public information available to indicate otherwise."
$ws.Range("L2").Value = "2025-03-12 08:14:33"
